$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12-21, columns C,D,E,F: replace "formula * 100" with raw fraction value
# and switch display to percentage format (0.00%) instead of manual *100 + "0.00".
$values = @{
    12 = @{ C = 0.98900002241134599; D = 0.98250001668929998; E = 0.98949998617172197; F = 0.99180001020431496 }
    13 = @{ C = 0.98369997739791804; D = 0.98820000886917103; E = 0.98570001125335605; F = 0.98830002546310403 }
    14 = @{ C = 0.97519999742507901; D = 0.98509997129440297; E = 0.99010002613067605; F = 0.98960000276565496 }
    15 = @{ C = 0.96380001306533802; D = 0.98360002040863004; E = 0.98780000209808305; F = 0.98600000143051103 }
    16 = @{ C = 0.95270001888275102; D = 0.98329997062683105; E = 0.98769998550414995; F = 0.99089998006820601 }
    17 = @{ C = 0.91990000009536699; D = 0.98540002107620195; E = 0.98790001869201605; F = 0.98900002241134599 }
    18 = @{ C = 0.75190001726150502; D = 0.98420000076293901; E = 0.98549997806548995; F = 0.98729997873306197 }
    19 = @{ C = 0.44710001349449102; D = 0.98309999704360895; E = 0.98619997501373202; F = 0.98530000448226895 }
    20 = @{ C = 0.100400000810623;   D = 0.97710001468658403; E = 0.98059999942779497; F = 0.98400002717971802 }
    21 = @{ C = 0.097400002181529999;D = 0.96310001611709495; E = 0.97490000724792403; F = 0.97790002822875899 }
}

foreach ($r in 12..21) {
    foreach ($col in 'C','D','E','F') {
        $cell = $ws.Range(($col + $r))
        $cell.Value = $values[$r][$col]
        $cell.NumberFormat = "0.00%"
    }
}

# Update the selected cell/active cell in the sheet view.
$ws.Range("K29").Select()
